# Regenerate save_data to use K instead of Strike# (column G, header "K").
# This recomputes the per-row K values (std/mean-derived s_vals calc) and
# writes the results back into column G for each data row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values keyed by worksheet row number (row 1 is the header row).
$kValues = @{
    2 = 1
    3 = 1
    4 = 3
    5 = 0
    6 = 1
    7 = 1
    8 = 1
    9 = 3
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 3
    15 = 1
    16 = 1
    17 = 0
    18 = 0
    19 = 1
    20 = 2
    21 = 0
    22 = 3
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 1
    30 = 2
    31 = 2
    32 = 0
    33 = 0
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 2
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 0
    51 = 0
    52 = 1
    53 = 1
    54 = 2
    55 = 1
    56 = 0
    58 = 1
    59 = 1
    60 = 1
    61 = 0
    62 = 1
    63 = 0
    64 = 0
    65 = 1
    66 = 1
    69 = 1
    70 = 2
    71 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
